$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('E7').Value = '  +1.34%  '
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('E11').Value = '  -2.18%  '
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('E14').Value = '  -2.98%  '
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('E23').Value = '  -4.02%  '
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('E29').Value = '  -4.78%  '
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E38').Value = '  -5.63%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E42').Value = '  -4.28%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E43').Value = '  -3.22%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E45').Value = '  -3.50%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E51').Value = '  -0.55%  '

# Price column (D) values must remain text; force text via quote-prefix then reset style
$ws.Range('D2').Value = "'26.406.46"
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'1.721.17"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Value = "'243.93"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Value = "'0.2611"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06196"
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Value = "'1.730.33"
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Value = "'15.40"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'4.531"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'0.5979"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'77.13"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'1.001"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'26.416.12"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Value = "'0.000007176"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'11.39"
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'1.953.45"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'4.490"
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = "'8.576"
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Value = "'5.176"
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Value = "'137.98"
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = "'15.23"
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'1.411"
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = "'107.13"
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Value = "'1.714"
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Value = "'3.957"
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Value = "'0.07942"
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Value = "'3.679"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = "'0.04529"
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Value = "'2.607"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Value = "'0.9947"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'0.6217"
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Value = "'0.9064"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'1.976"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Value = "'2.399"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'1.001"
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'0.01484"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'100.31"
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = "'5.398"
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = "'0.3847"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'6.736"
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'0.1149"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'0.05358"
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = "'30.12"
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'7.695"
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'1.249"
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = "'51.06"
$ws.Range('D51').Style = 'Normal'
